# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header style (copied from the adjacent "sum" header in G1) and filling in
# the per-row save indicator values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same bold/centered/bordered style as the other headers.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row "Save" flag values (rows 2 through 13).
$values = @(1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
